# Insert a new data row into the "Hortaliza, Macroferia Regional de Talca - Zanahoria"
# sheet right before the existing row 422, shifting it (and every row below it)
# down by one, then populate the newly-opened row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 422 (and everything after it) down by one row.
$ws.Rows.Item(422).Insert()

# Fill in the new record in the now-empty row 422.
$ws.Range("A422").Value = 5
$ws.Range("B422").Value = "Macroferia Regional de Talca"
$ws.Range("C422").Value = "Maule"
$ws.Range("D422").Value = 44984
$ws.Range("E422").Value = 7
$ws.Range("F422").Value = 100114013
$ws.Range("G422").Value = "Zanahoria"
$ws.Range("H422").Value = "Sin especificar"
$ws.Range("I422").Value = "Primera"
$ws.Range("J422").Value = 500
$ws.Range("K422").Value = 8000
$ws.Range("L422").Value = 8000
$ws.Range("M422").Value = 8000
$ws.Range("N422").Value = "$/saco 20 kilos"
$ws.Range("O422").Value = "Región de La Araucanía"
$ws.Range("P422").Value = 400
$ws.Range("Q422").Value = 20
$ws.Range("R422").Value = "Hortaliza"
